$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column E (pt_min -> pt_max will go here)
$ws.Columns("E:E").Insert()

$ws.Range("E1").Value = "pt_max"
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 5).Value = 50
}

$ws.Range("E17").Select()
